$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "ALC"; Cell = "H6"; Value = 308.66666 },
    @{ Sheet = "ALC"; Cell = "I6"; Value = 130 },
    @{ Sheet = "ALC"; Cell = "J6"; Value = 398 },
    @{ Sheet = "ALC"; Cell = "K6"; Value = 390 },
    @{ Sheet = "ALC"; Cell = "L6"; Value = 1194 },
    @{ Sheet = "ALC"; Cell = "M6"; Value = -278 },
    @{ Sheet = "ALC"; Cell = "N6"; Value = -1418 },
    @{ Sheet = "ALC"; Cell = "H98"; Value = 1484.5333 },
    @{ Sheet = "ALC"; Cell = "I98"; Value = 831.8 },
    @{ Sheet = "ALC"; Cell = "K98"; Value = 831.8 },
    @{ Sheet = "ALC"; Cell = "M98"; Value = 666.2 },
    @{ Sheet = "ALC"; Cell = "H122"; Value = 1484.5333 },
    @{ Sheet = "ALC"; Cell = "I122"; Value = 831.8 },
    @{ Sheet = "ALC"; Cell = "K122"; Value = 2495.4 },
    @{ Sheet = "ALC"; Cell = "M122"; Value = -45.39999999999964 },
    @{ Sheet = "ALC"; Cell = "H137"; Value = 4168069 },
    @{ Sheet = "ALC"; Cell = "I137"; Value = 2084628.5 },
    @{ Sheet = "ALC"; Cell = "K137"; Value = 6253885.5 },
    @{ Sheet = "ALC"; Cell = "M137"; Value = -6251335.5 },
    @{ Sheet = "ARM"; Cell = "H32"; Value = 20714.725 },
    @{ Sheet = "ARM"; Cell = "I32"; Value = 4390.4395 },
    @{ Sheet = "ARM"; Cell = "J32"; Value = 128455 },
    @{ Sheet = "ARM"; Cell = "K32"; Value = 4390.4395 },
    @{ Sheet = "ARM"; Cell = "L32"; Value = 128455 },
    @{ Sheet = "ARM"; Cell = "M32"; Value = -4103.4395 },
    @{ Sheet = "ARM"; Cell = "N32"; Value = -129029 },
    @{ Sheet = "ARM"; Cell = "H74"; Value = 979.3714 },
    @{ Sheet = "ARM"; Cell = "I74"; Value = 973.4838999999999 },
    @{ Sheet = "ARM"; Cell = "J74"; Value = 1025 },
    @{ Sheet = "ARM"; Cell = "K74"; Value = 973.4838999999999 },
    @{ Sheet = "ARM"; Cell = "L74"; Value = 1025 },
    @{ Sheet = "ARM"; Cell = "M74"; Value = -99.48389999999995 },
    @{ Sheet = "ARM"; Cell = "N74"; Value = -2773 },
    @{ Sheet = "ARM"; Cell = "H77"; Value = 979.3714 },
    @{ Sheet = "ARM"; Cell = "I77"; Value = 973.4838999999999 },
    @{ Sheet = "ARM"; Cell = "J77"; Value = 1025 },
    @{ Sheet = "ARM"; Cell = "K77"; Value = 4867.4195 },
    @{ Sheet = "ARM"; Cell = "L77"; Value = 5125 },
    @{ Sheet = "ARM"; Cell = "M77"; Value = -499.4195 },
    @{ Sheet = "ARM"; Cell = "N77"; Value = -13861 },
    @{ Sheet = "ARM"; Cell = "H97"; Value = 2017.5 },
    @{ Sheet = "ARM"; Cell = "I97"; Value = 1076.4286 },
    @{ Sheet = "ARM"; Cell = "J97"; Value = 2840.9375 },
    @{ Sheet = "ARM"; Cell = "K97"; Value = 1076.4286 },
    @{ Sheet = "ARM"; Cell = "L97"; Value = 2840.9375 },
    @{ Sheet = "ARM"; Cell = "M97"; Value = -580.4286 },
    @{ Sheet = "ARM"; Cell = "N97"; Value = -3832.9375 },
    @{ Sheet = "BSM"; Cell = "H94"; Value = 799.9474 },
    @{ Sheet = "BSM"; Cell = "I94"; Value = 800.75 },
    @{ Sheet = "BSM"; Cell = "J94"; Value = 798.5714 },
    @{ Sheet = "BSM"; Cell = "K94"; Value = 800.75 },
    @{ Sheet = "BSM"; Cell = "L94"; Value = 798.5714 },
    @{ Sheet = "BSM"; Cell = "M94"; Value = -349.75 },
    @{ Sheet = "BSM"; Cell = "N94"; Value = -1700.5714 },
    @{ Sheet = "BSM"; Cell = "H112"; Value = 40660 },
    @{ Sheet = "BSM"; Cell = "J112"; Value = 40660 },
    @{ Sheet = "BSM"; Cell = "L112"; Value = 40660 },
    @{ Sheet = "BSM"; Cell = "N112"; Value = -43614 },
    @{ Sheet = "BSM"; Cell = "H134"; Value = 252385.33 },
    @{ Sheet = "BSM"; Cell = "I134"; Value = 376828 },
    @{ Sheet = "BSM"; Cell = "J134"; Value = 3500 },
    @{ Sheet = "BSM"; Cell = "K134"; Value = 1130484 },
    @{ Sheet = "BSM"; Cell = "L134"; Value = 10500 },
    @{ Sheet = "BSM"; Cell = "M134"; Value = -1127949 },
    @{ Sheet = "BSM"; Cell = "N134"; Value = -15570 },
    @{ Sheet = "CRP"; Cell = "H31"; Value = 15153648 },
    @{ Sheet = "CRP"; Cell = "I31"; Value = 1814.125 },
    @{ Sheet = "CRP"; Cell = "K31"; Value = 1814.125 },
    @{ Sheet = "CRP"; Cell = "M31"; Value = -1519.125 },
    @{ Sheet = "CRP"; Cell = "H34"; Value = 15153648 },
    @{ Sheet = "CRP"; Cell = "I34"; Value = 1814.125 },
    @{ Sheet = "CRP"; Cell = "K34"; Value = 1814.125 },
    @{ Sheet = "CRP"; Cell = "M34"; Value = -1612.125 },
    @{ Sheet = "CRP"; Cell = "H58"; Value = 826.3509 },
    @{ Sheet = "CRP"; Cell = "I58"; Value = 845.6889 },
    @{ Sheet = "CRP"; Cell = "K58"; Value = 845.6889 },
    @{ Sheet = "CRP"; Cell = "M58"; Value = -642.6889 },
    @{ Sheet = "CRP"; Cell = "H134"; Value = 6673868.5 },
    @{ Sheet = "CRP"; Cell = "I134"; Value = 8433.666999999999 },
    @{ Sheet = "CRP"; Cell = "K134"; Value = 25301.001 },
    @{ Sheet = "CRP"; Cell = "M134"; Value = -22766.001 },
    @{ Sheet = "CRP"; Cell = "H136"; Value = 826.3509 },
    @{ Sheet = "CRP"; Cell = "I136"; Value = 845.6889 },
    @{ Sheet = "CRP"; Cell = "K136"; Value = 2537.0667 },
    @{ Sheet = "CRP"; Cell = "M136"; Value = 12.93330000000014 },
    @{ Sheet = "CUL"; Cell = "H37"; Value = 73333.336 },
    @{ Sheet = "CUL"; Cell = "J37"; Value = 73333.336 },
    @{ Sheet = "CUL"; Cell = "L37"; Value = 220000.008 },
    @{ Sheet = "CUL"; Cell = "N37"; Value = -220224.008 },
    @{ Sheet = "CUL"; Cell = "H131"; Value = 863.24 },
    @{ Sheet = "CUL"; Cell = "I131"; Value = 489.07693 },
    @{ Sheet = "CUL"; Cell = "J131"; Value = 994.7027 },
    @{ Sheet = "CUL"; Cell = "K131"; Value = 1467.23079 },
    @{ Sheet = "CUL"; Cell = "L131"; Value = 2984.1081 },
    @{ Sheet = "CUL"; Cell = "M131"; Value = 3572.76921 },
    @{ Sheet = "CUL"; Cell = "N131"; Value = -13064.1081 },
    @{ Sheet = "CUL"; Cell = "H132"; Value = 693990 },
    @{ Sheet = "CUL"; Cell = "I132"; Value = 940772.1 },
    @{ Sheet = "CUL"; Cell = "K132"; Value = 8466948.9 },
    @{ Sheet = "CUL"; Cell = "M132"; Value = -8464418.9 },
    @{ Sheet = "CUL"; Cell = "H137"; Value = 2069.6924 },
    @{ Sheet = "CUL"; Cell = "J137"; Value = 2261 },
    @{ Sheet = "CUL"; Cell = "L137"; Value = 6783 },
    @{ Sheet = "CUL"; Cell = "N137"; Value = -16983 },
    @{ Sheet = "GSM"; Cell = "H97"; Value = 1589.8 },
    @{ Sheet = "GSM"; Cell = "I97"; Value = 1707.125 },
    @{ Sheet = "GSM"; Cell = "J97"; Value = 1455.7142 },
    @{ Sheet = "GSM"; Cell = "K97"; Value = 1707.125 },
    @{ Sheet = "GSM"; Cell = "L97"; Value = 1455.7142 },
    @{ Sheet = "GSM"; Cell = "M97"; Value = -1211.125 },
    @{ Sheet = "GSM"; Cell = "N97"; Value = -2447.7142 },
    @{ Sheet = "GSM"; Cell = "H122"; Value = 2266.5264 },
    @{ Sheet = "GSM"; Cell = "I122"; Value = 2263.6667 },
    @{ Sheet = "GSM"; Cell = "J122"; Value = 2271.4285 },
    @{ Sheet = "GSM"; Cell = "K122"; Value = 6791.000100000001 },
    @{ Sheet = "GSM"; Cell = "L122"; Value = 6814.2855 },
    @{ Sheet = "GSM"; Cell = "M122"; Value = -4341.000100000001 },
    @{ Sheet = "GSM"; Cell = "N122"; Value = -11714.2855 },
    @{ Sheet = "GSM"; Cell = "H126"; Value = 12316.25 },
    @{ Sheet = "GSM"; Cell = "I126"; Value = 2156.889 },
    @{ Sheet = "GSM"; Cell = "J126"; Value = 18411.867 },
    @{ Sheet = "GSM"; Cell = "K126"; Value = 6470.667 },
    @{ Sheet = "GSM"; Cell = "L126"; Value = 55235.601 },
    @{ Sheet = "GSM"; Cell = "M126"; Value = -4000.667 },
    @{ Sheet = "GSM"; Cell = "N126"; Value = -60175.601 },
    @{ Sheet = "GSM"; Cell = "H132"; Value = 4999.3335 },
    @{ Sheet = "GSM"; Cell = "I132"; Value = 5000 },
    @{ Sheet = "GSM"; Cell = "J132"; Value = 4999 },
    @{ Sheet = "GSM"; Cell = "K132"; Value = 15000 },
    @{ Sheet = "GSM"; Cell = "L132"; Value = 14997 },
    @{ Sheet = "GSM"; Cell = "M132"; Value = -12470 },
    @{ Sheet = "GSM"; Cell = "N132"; Value = -20057 },
    @{ Sheet = "LTW"; Cell = "H40"; Value = 2199.2222 },
    @{ Sheet = "LTW"; Cell = "I40"; Value = 1818.6 },
    @{ Sheet = "LTW"; Cell = "J40"; Value = 2675 },
    @{ Sheet = "LTW"; Cell = "K40"; Value = 1818.6 },
    @{ Sheet = "LTW"; Cell = "L40"; Value = 2675 },
    @{ Sheet = "LTW"; Cell = "M40"; Value = -1682.6 },
    @{ Sheet = "LTW"; Cell = "N40"; Value = -2947 },
    @{ Sheet = "LTW"; Cell = "H61"; Value = 2446.1538 },
    @{ Sheet = "LTW"; Cell = "I61"; Value = 1950 },
    @{ Sheet = "LTW"; Cell = "J61"; Value = 4100 },
    @{ Sheet = "LTW"; Cell = "K61"; Value = 1950 },
    @{ Sheet = "LTW"; Cell = "L61"; Value = 4100 },
    @{ Sheet = "LTW"; Cell = "M61"; Value = -1748 },
    @{ Sheet = "LTW"; Cell = "N61"; Value = -4504 },
    @{ Sheet = "LTW"; Cell = "H68"; Value = 2460 },
    @{ Sheet = "LTW"; Cell = "I68"; Value = 1433.3334 },
    @{ Sheet = "LTW"; Cell = "J68"; Value = 4000 },
    @{ Sheet = "LTW"; Cell = "K68"; Value = 1433.3334 },
    @{ Sheet = "LTW"; Cell = "L68"; Value = 4000 },
    @{ Sheet = "LTW"; Cell = "M68"; Value = -684.3334 },
    @{ Sheet = "LTW"; Cell = "N68"; Value = -5498 },
    @{ Sheet = "LTW"; Cell = "H71"; Value = 2460 },
    @{ Sheet = "LTW"; Cell = "I71"; Value = 1433.3334 },
    @{ Sheet = "LTW"; Cell = "J71"; Value = 4000 },
    @{ Sheet = "LTW"; Cell = "K71"; Value = 7166.666999999999 },
    @{ Sheet = "LTW"; Cell = "L71"; Value = 20000 },
    @{ Sheet = "LTW"; Cell = "M71"; Value = -3422.666999999999 },
    @{ Sheet = "LTW"; Cell = "N71"; Value = -27488 },
    @{ Sheet = "LTW"; Cell = "H113"; Value = 2446.1538 },
    @{ Sheet = "LTW"; Cell = "I113"; Value = 1950 },
    @{ Sheet = "LTW"; Cell = "J113"; Value = 4100 },
    @{ Sheet = "LTW"; Cell = "K113"; Value = 1950 },
    @{ Sheet = "LTW"; Cell = "L113"; Value = 4100 },
    @{ Sheet = "LTW"; Cell = "M113"; Value = 220 },
    @{ Sheet = "LTW"; Cell = "N113"; Value = -8440 },
    @{ Sheet = "LTW"; Cell = "H136"; Value = 1987.2084 },
    @{ Sheet = "LTW"; Cell = "I136"; Value = 1769.7 },
    @{ Sheet = "LTW"; Cell = "J136"; Value = 3074.75 },
    @{ Sheet = "LTW"; Cell = "K136"; Value = 5309.1 },
    @{ Sheet = "LTW"; Cell = "L136"; Value = 9224.25 },
    @{ Sheet = "LTW"; Cell = "M136"; Value = -2759.1 },
    @{ Sheet = "LTW"; Cell = "N136"; Value = -14324.25 },
    @{ Sheet = "WVR"; Cell = "H62"; Value = 3789.111 },
    @{ Sheet = "WVR"; Cell = "I62"; Value = 3643.1428 },
    @{ Sheet = "WVR"; Cell = "J62"; Value = 4300 },
    @{ Sheet = "WVR"; Cell = "K62"; Value = 3643.1428 },
    @{ Sheet = "WVR"; Cell = "L62"; Value = 4300 },
    @{ Sheet = "WVR"; Cell = "M62"; Value = -3019.1428 },
    @{ Sheet = "WVR"; Cell = "N62"; Value = -5548 },
    @{ Sheet = "WVR"; Cell = "H65"; Value = 3789.111 },
    @{ Sheet = "WVR"; Cell = "I65"; Value = 3643.1428 },
    @{ Sheet = "WVR"; Cell = "J65"; Value = 4300 },
    @{ Sheet = "WVR"; Cell = "K65"; Value = 18215.714 },
    @{ Sheet = "WVR"; Cell = "L65"; Value = 21500 },
    @{ Sheet = "WVR"; Cell = "M65"; Value = -15095.714 },
    @{ Sheet = "WVR"; Cell = "N65"; Value = -27740 },
    @{ Sheet = "WVR"; Cell = "H96"; Value = 1666.6666 },
    @{ Sheet = "WVR"; Cell = "I96"; Value = 1500 },
    @{ Sheet = "WVR"; Cell = "J96"; Value = 2000 },
    @{ Sheet = "WVR"; Cell = "K96"; Value = 1500 },
    @{ Sheet = "WVR"; Cell = "L96"; Value = 2000 },
    @{ Sheet = "WVR"; Cell = "M96"; Value = -127 },
    @{ Sheet = "WVR"; Cell = "N96"; Value = -4746 },
    @{ Sheet = "WVR"; Cell = "H132"; Value = 2329.9092 },
    @{ Sheet = "WVR"; Cell = "I132"; Value = 2492.1714 },
    @{ Sheet = "WVR"; Cell = "J132"; Value = 1698.8889 },
    @{ Sheet = "WVR"; Cell = "K132"; Value = 7476.514200000001 },
    @{ Sheet = "WVR"; Cell = "L132"; Value = 5096.6667 },
    @{ Sheet = "WVR"; Cell = "M132"; Value = -4946.514200000001 },
    @{ Sheet = "WVR"; Cell = "N132"; Value = -10156.6667 },
    @{ Sheet = "WVR"; Cell = "H136"; Value = 1773.4546 },
    @{ Sheet = "WVR"; Cell = "I136"; Value = 1924 },
    @{ Sheet = "WVR"; Cell = "J136"; Value = 1214.2858 },
    @{ Sheet = "WVR"; Cell = "K136"; Value = 5772 },
    @{ Sheet = "WVR"; Cell = "L136"; Value = 3642.8574 },
    @{ Sheet = "WVR"; Cell = "M136"; Value = -3222 },
    @{ Sheet = "WVR"; Cell = "N136"; Value = -8742.857400000001 }
)

$wsCache = @{}
foreach ($c in $changes) {
    if (-not $wsCache.ContainsKey($c.Sheet)) {
        $wsCache[$c.Sheet] = $wb.Worksheets.Item($c.Sheet)
    }
    $ws = $wsCache[$c.Sheet]
    $ws.Range($c.Cell).Value = $c.Value
}